$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53: count of crashes
$ws.Range("F53").Value = "N. crash:"
$ws.Range("G53").Formula = "=COUNTIF(G2:G51,""True"")"

# Row 54: mean of columns I:P
$ws.Range("F54").Value = "Mean:"
$ws.Range("I54:P54").Formula = "=AVERAGE(I2:I51)"

# Row 55: standard deviation of columns I:P
$ws.Range("F55").Value = "Standard Deviation:"
$ws.Range("I55:P55").Formula = "=STDEV.S(I2:I51)"

# Update the view: scroll window and selection to match the saved state
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("Z38").Select()
